# MitsosBarton2006Ex314 - Strong_Stationary generator (alpha non zero)
# nuevos experimentos no convexos
#
# Re-writes a handful of text-valued "label" cells (expressions / numeric
# strings stored as text) on several sheets, plus one genuinely-numeric
# cell on Vector_Alpha.
#
# Cells such as B2 on "Restricciones_del_lider" look numeric ("1.1") but
# are stored as shared-string TEXT in the workbook, not as numbers. Simply
# assigning a numeric-looking string to .Value lets Excel auto-convert it
# to a number cell, which would not match the original text-cell layout.
# To keep them as text we flip NumberFormat to "@" (Text) before writing
# the value, then reset the cell Style back to "Normal" so no stray
# number-format style lingers on the cell once we're done.
#
# NOTE: worksheets are addressed by their 1-based index rather than by
# name. The workbook has two sheets whose names differ only by case
# ("Vector_bf" / "Vector_BF") and name-based lookup resolves both to the
# same (first) sheet, so index-based access is used everywhere to stay
# unambiguous.
#   1 Funciones_Objetivo
#   2 Restricciones_del_lider
#   3 Restricciones_del_follower
#   4 Punto_modificado
#   5 Vector_bf
#   6 Vector_BF
#   7 Vector_Alpha

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param(
        $Worksheet,
        [string]$CellRef,
        [string]$Text
    )
    $cell = $Worksheet.Range($CellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.Style = "Normal"
}

# --- Restricciones_del_lider (sheet 2) ------------------------------------
$ws = $wb.Worksheets.Item(2)
Set-TextValue $ws "A2" "-0.9 + x"
Set-TextValue $ws "B2" "-0.09999999999999998"
Set-TextValue $ws "D2" "0.42"
Set-TextValue $ws "A3" "0.8999999999999999 - x"
Set-TextValue $ws "B3" "-1.9"
Set-TextValue $ws "D3" "0.02"

# --- Restricciones_del_follower (sheet 3) ---------------------------------
$ws = $wb.Worksheets.Item(3)
Set-TextValue $ws "A2" "-24.666666666666668 + 9.135802469135802y"
Set-TextValue $ws "B2" "23.666666666666668"
Set-TextValue $ws "D2" "0.29"
Set-TextValue $ws "E2" "0.8"
Set-TextValue $ws "F2" "7.4"
Set-TextValue $ws "A3" "0"
Set-TextValue $ws "B3" "-1"
Set-TextValue $ws "D3" "0.52"
Set-TextValue $ws "E3" "0.2"
Set-TextValue $ws "F3" "6.0"

# --- Punto_modificado (sheet 4) -------------------------------------------
$ws = $wb.Worksheets.Item(4)
Set-TextValue $ws "A2" "0.9"
Set-TextValue $ws "B2" "2.7"

# --- Vector_bf (sheet 5) ---------------------------------------------------
$ws = $wb.Worksheets.Item(5)
Set-TextValue $ws "A2" "-9.039382716049383"

# --- Vector_BF (sheet 6) ----------------------------------------------------
$ws = $wb.Worksheets.Item(6)
Set-TextValue $ws "A2" "-0.8899999999999999"
Set-TextValue $ws "A3" "-17.082641975308643"

# --- Vector_Alpha (sheet 7, genuine numeric cell) ---------------------------
$ws = $wb.Worksheets.Item(7)
$ws.Range("A2").Value = 0.81

Write-Output "edits applied"
